$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 122, shifting rows 122:173 down to 123:174
$ws.Rows.Item(122).Insert()

# Populate the new row 122 with the new data.
$ws.Cells.Item(122, 1).Value = 7
$ws.Cells.Item(122, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value = "Ñuble"
$ws.Cells.Item(122, 4).Value = 45141
$ws.Cells.Item(122, 4).NumberFormat = $ws.Cells.Item(123, 4).NumberFormat
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(122, 6).Value = "Fruta"
$ws.Cells.Item(122, 7).Value = 100108
$ws.Cells.Item(122, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(122, 9).Value = 100108002
$ws.Cells.Item(122, 10).Value = "Mango"
$ws.Cells.Item(122, 11).Value = "Sin especificar"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 50
$ws.Cells.Item(122, 14).Value = 8000
$ws.Cells.Item(122, 15).Value = 8000
$ws.Cells.Item(122, 16).Value = 8000
$ws.Cells.Item(122, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(122, 18).Value = "Brasil"
$ws.Cells.Item(122, 19).Value = 2000
$ws.Cells.Item(122, 20).Value = 4
